$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 6000
$ws.Range("J8").Value = 6000
$ws.Range("L8").Value = 18000
$ws.Range("N8").Value = -18278
$ws.Range("H80").Value = 4292.3125
$ws.Range("I80").Value = 866.5
$ws.Range("K80").Value = 2599.5
$ws.Range("M80").Value = -1601.5
$ws.Range("H82").Value = 1558.3334
$ws.Range("I82").Value = 1558.3334
$ws.Range("K82").Value = 4675.0002
$ws.Range("M82").Value = -4269.0002
$ws.Range("H83").Value = 4292.3125
$ws.Range("I83").Value = 866.5
$ws.Range("K83").Value = 7798.5
$ws.Range("M83").Value = -2806.5
$ws.Range("H85").Value = 1558.3334
$ws.Range("I85").Value = 1558.3334
$ws.Range("K85").Value = 4675.0002
$ws.Range("M85").Value = -3271.0002
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H116").Value = 2641
$ws.Range("I116").Value = 1735
$ws.Range("K116").Value = 1735
$ws.Range("M116").Value = 1707
$ws.Range("H118").Value = 373.57144
$ws.Range("I118").Value = 243
$ws.Range("J118").Value = 700
$ws.Range("K118").Value = 729
$ws.Range("L118").Value = 2100
$ws.Range("M118").Value = 928
$ws.Range("N118").Value = -5414
$ws.Range("H132").Value = 1843
$ws.Range("I132").Value = 1576.6923
$ws.Range("K132").Value = 4730.0769
$ws.Range("M132").Value = -2200.0769
$ws.Range("H138").Value = 3866.275
$ws.Range("I138").Value = 1978.5714
$ws.Range("J138").Value = 4882.731
$ws.Range("K138").Value = 5935.7142
$ws.Range("L138").Value = 14648.193
$ws.Range("M138").Value = -795.7142000000003
$ws.Range("N138").Value = -24928.193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1931.6666
$ws.Range("I2").Value = 2002.0952
$ws.Range("J2").Value = 1685.1666
$ws.Range("K2").Value = 2002.0952
$ws.Range("L2").Value = 1685.1666
$ws.Range("M2").Value = -1889.0952
$ws.Range("N2").Value = -1911.1666
$ws.Range("H34").Value = 36957
$ws.Range("I34").Value = 8000
$ws.Range("J34").Value = 46609.332
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 46609.332
$ws.Range("M34").Value = -7729
$ws.Range("N34").Value = -47151.332
$ws.Range("H45").Value = 1906.2609
$ws.Range("I45").Value = 1758.7368
$ws.Range("J45").Value = 2607
$ws.Range("K45").Value = 1758.7368
$ws.Range("L45").Value = 2607
$ws.Range("M45").Value = -1381.7368
$ws.Range("N45").Value = -3361
$ws.Range("H74").Value = 204663.28
$ws.Range("I74").Value = 223845.31
$ws.Range("J74").Value = 51207
$ws.Range("K74").Value = 223845.31
$ws.Range("L74").Value = 51207
$ws.Range("M74").Value = -222971.31
$ws.Range("N74").Value = -52955
$ws.Range("H77").Value = 204663.28
$ws.Range("I77").Value = 223845.31
$ws.Range("J77").Value = 51207
$ws.Range("K77").Value = 1119226.55
$ws.Range("L77").Value = 256035
$ws.Range("M77").Value = -1114858.55
$ws.Range("N77").Value = -264771
$ws.Range("H116").Value = 1931.6666
$ws.Range("I116").Value = 2002.0952
$ws.Range("J116").Value = 1685.1666
$ws.Range("K116").Value = 2002.0952
$ws.Range("L116").Value = 1685.1666
$ws.Range("M116").Value = 291.9048
$ws.Range("N116").Value = -6273.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1931.6666
$ws.Range("I3").Value = 2002.0952
$ws.Range("J3").Value = 1685.1666
$ws.Range("K3").Value = 2002.0952
$ws.Range("L3").Value = 1685.1666
$ws.Range("M3").Value = -1888.0952
$ws.Range("N3").Value = -1913.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1414.8334
$ws.Range("I16").Value = 1395.375
$ws.Range("J16").Value = 1453.75
$ws.Range("K16").Value = 1395.375
$ws.Range("L16").Value = 1453.75
$ws.Range("M16").Value = -1108.375
$ws.Range("N16").Value = -2027.75
$ws.Range("H107").Value = 1409.2727
$ws.Range("I107").Value = 1390.2
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 1390.2
$ws.Range("L107").Value = 1600
$ws.Range("M107").Value = 529.8
$ws.Range("N107").Value = -5440
$ws.Range("H113").Value = 1414.8334
$ws.Range("I113").Value = 1395.375
$ws.Range("J113").Value = 1453.75
$ws.Range("K113").Value = 1395.375
$ws.Range("L113").Value = 1453.75
$ws.Range("M113").Value = 774.625
$ws.Range("N113").Value = -5793.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12611.038
$ws.Range("I131").Value = 459.89474
$ws.Range("J131").Value = 44093.547
$ws.Range("K131").Value = 1379.68422
$ws.Range("L131").Value = 132280.641
$ws.Range("M131").Value = 3660.31578
$ws.Range("N131").Value = -142360.641
$ws.Range("H136").Value = 2257.6924
$ws.Range("I136").Value = 981.93335
$ws.Range("J136").Value = 3997.3635
$ws.Range("K136").Value = 2945.80005
$ws.Range("L136").Value = 11992.0905
$ws.Range("M136").Value = 2154.19995
$ws.Range("N136").Value = -22192.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2092.2083
$ws.Range("I113").Value = 2111.0527
$ws.Range("J113").Value = 2020.6
$ws.Range("K113").Value = 2111.0527
$ws.Range("L113").Value = 2020.6
$ws.Range("M113").Value = 58.94729999999981
$ws.Range("N113").Value = -6360.6
$ws.Range("H126").Value = 3010.7
$ws.Range("I126").Value = 1950.1
$ws.Range("J126").Value = 4071.3
$ws.Range("K126").Value = 5850.299999999999
$ws.Range("L126").Value = 12213.9
$ws.Range("M126").Value = -3380.299999999999
$ws.Range("N126").Value = -17153.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 817595.7
$ws.Range("I61").Value = 959560.1
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 959560.1
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -959358.1
$ws.Range("N61").Value = -1704
$ws.Range("H82").Value = 2299.0908
$ws.Range("I82").Value = 1348.75
$ws.Range("J82").Value = 4833.3335
$ws.Range("K82").Value = 1348.75
$ws.Range("L82").Value = 4833.3335
$ws.Range("M82").Value = -987.75
$ws.Range("N82").Value = -5555.3335
$ws.Range("H85").Value = 2299.0908
$ws.Range("I85").Value = 1348.75
$ws.Range("J85").Value = 4833.3335
$ws.Range("K85").Value = 1348.75
$ws.Range("L85").Value = 4833.3335
$ws.Range("M85").Value = -100.75
$ws.Range("N85").Value = -7329.3335
$ws.Range("H113").Value = 817595.7
$ws.Range("I113").Value = 959560.1
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 959560.1
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = -957390.1
$ws.Range("N113").Value = -5640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 77503
$ws.Range("I3").Value = 10000
$ws.Range("K3").Value = 10000
$ws.Range("M3").Value = -9886
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H113").Value = 7602.533
$ws.Range("I113").Value = 813.6
$ws.Range("J113").Value = 21180.4
$ws.Range("K113").Value = 2440.8
$ws.Range("L113").Value = 63541.2
$ws.Range("M113").Value = -270.8000000000002
$ws.Range("N113").Value = -67881.20000000001

Write-Host "edit.ps1 applied 187 cell changes across 8 sheets"
